$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C1").Value = "Period [s] (5 Hz)"
$ws.Range("D1").Value = "Period [s] (6 Hz)"

$ws.Activate()
$ws.Range("B1").Select()
